# Generate Report for Handback
# For the file "2465d413-ea25-4cde-946a-75d0f60ed71a.md", the handback transform
# failed (the generated handback file name did not match the handoff file name).
# Update the Overview sheet and the per-locale (zh-cn / de-de) sheets to reflect
# the new status and the corresponding error detail message.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E7").Value = "Handback transform failed"
$overview.Range("F7").Value = "Handback transform failed"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C7").Value = "Handback transform failed"
$zhcn.Range("P7").Value = "Handback file name: exe1xz1o.ucl is different with handoff file name: 2465d413-ea25-4cde-946a-75d0f60ed71a.149bd73d9689c2ce6138abf6956767c53ae7dc60.zh-cn."

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C7").Value = "Handback transform failed"
$dede.Range("P7").Value = "Handback file name: exe1xz1o.ucl is different with handoff file name: 2465d413-ea25-4cde-946a-75d0f60ed71a.149bd73d9689c2ce6138abf6956767c53ae7dc60.de-de."
